$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old data column (A1:A9) that held the protocol numbers.
# ---------------------------------------------------------------------------
$ws.Range("A1:A9").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. Touch F7 so it carries a (non-default) cell style, matching the new
#    single-cell sheet whose only populated cell is F7 with style index 1.
# ---------------------------------------------------------------------------
$ws.Range("F7").Font.Name = "Calibri"

# Move the selection / active cell to F7 (was F9 before).
$ws.Range("F7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update header text: "O.S. Corretiva" -> "O.S. Corretiva/Preventiva"
# ---------------------------------------------------------------------------
$ws.PageSetup.CenterHeader = "Protocolo de entrega`nO.S. Corretiva/Preventiva"

# ---------------------------------------------------------------------------
# 4. Rebuild the footer: right section keeps the signature line, left
#    section becomes a coloured/sized "Data" line plus a classification
#    line, separated by a carriage return.
# ---------------------------------------------------------------------------
$cr = [char]13
$leftFooter = '&"Calibri"&11&K000000      Data: _____/_____/________' + $cr + '&1#&"Calibri"&10&K0000FFClassificação: Interno'
$rightFooter = "Assinatura do(a) responsavel: _____________________________"

$ws.PageSetup.LeftFooter = $leftFooter
$ws.PageSetup.CenterFooter = ""
$ws.PageSetup.RightFooter = $rightFooter
